$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# ---------------------------------------------------------------------------
# New "Friday" day-separator row (846) - copy the look of an existing
# day-separator row (row 5) and stamp the day name into column B.
# ---------------------------------------------------------------------------
$ws.Range("A5:F5").Copy()
$ws.Range("A846:F846").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B846").Value = "Friday"

# ---------------------------------------------------------------------------
# Helper data for the new log rows (847-859). Columns follow the existing
# sheet layout: A=Task type, B=Date, C=Time, D=Building, E=Room, F=Notes.
# ---------------------------------------------------------------------------
$newRows = @(
    @{ Row = 847; A = "Setup Mic"; C = "1715"; D = "LAS"; E = "C"; F = "Take cart with mixer, 2 wireless mics and 2 mic stands from Lassonde 1011 storeroom (across from Lassonde A). Go to Lassonde C classroom (class starts at 5:30 pm but be there early in case previous class ends early). "; Height = 60 },
    @{ Row = 848; A = "Other"; C = "1715"; D = "LAS"; E = "C"; F = 'Log in as 5065*0 on touchscreen. (First level bar is your wireless handheld mic volume). Plug in mic cable from output of mixer to mic input on podium (XLR jack just above VHS machine in podium). Ramp up volume a bit on "Microphone 2" on touchscreen to medium volume to get level.'; Height = 75 },
    @{ Row = 849; A = "Other"; C = "1715"; D = "LAS"; E = "C"; F = 'Plug in power cord from cart on to power outlet on left side of podium (to left of document camera). Turn on mixer. Turn on wireless microphone receivers on cart (NOTE: DO NOT PRESS "SYNC" BUTTON" - POWER BUTTON IS FIRST BUTTON TO THE RIGHT ON RECEIVER). '; Height = 75 },
    @{ Row = 850; A = "Other"; C = "1715"; D = "LAS"; E = "C"; F = 'Press "MUTE" button on wireless mics to turn on mics.'; Height = 0 },
    @{ Row = 851; A = "Other"; C = "1715"; D = "LAS"; E = "C"; F = "Once volumes are set, place one mic stand with mic halfway up aisle on right and one mic stand with mic halfway up aisle on left. Demo volume controls to prof. and demo PC. Leave microphone bags with milk carton on cart in room. PLEASE FIND OUT END TIME OF CLASS FROM PROF. AND TELL MASI AS MICROPHONES ARE EXPENSIVE. TELL PROF. TO STAY WITH MICS UNTIL THEY ARE PICKED UP. TELL HIM TO CALL ext 55800   WHEN DONE (use phone in classroom)."; Height = 120 },
    @{ Row = 852; A = "AV Shutdown"; C = "1730"; D = "CLH"; E = "L"; F = "PLEASE MAKE SURE CRESTRON GETS LOGGED OFF. WE ARE HAVING PROBLEMS WITH THIS ROOM WHEN IT DOESN'T GET LOGGED OFF."; Height = 45 },
    @{ Row = 853; A = "AV Shutdown"; C = "1900"; D = "LSB"; E = "103"; F = "Make sure neck mic goes back to drawer and log off touchscreen."; Height = 0 },
    @{ Row = 854; A = "AV Shutdown"; C = "1900"; D = "LSB"; E = "106"; F = "Make sure neck mic goes back to drawer and log off touchscreen."; Height = 0 },
    @{ Row = 855; A = "AV Shutdown"; C = "1900"; D = "LSB"; E = "107"; F = "Make sure neck mic goes back to drawer and log off touchscreen."; Height = 0 },
    @{ Row = 856; A = "Pickup Mic"; C = "1820"; D = "LAS"; E = "C"; F = "Pick up 2 wireless mics on stands with cart. Move all equipment on cart - cart has 2 wireless mic receivers and mixer and mic cables. Pick up 2 mic stands - return all equipment to Lassonde 1011 storeroom (across the hall from Lassonde A). PLEASE PUT 2 WIRELESS MICS IN BAGS PROVIDED IN MILK CARTON ON CART. Very expensive mics - please go early and treat mics with care."; Height = 90 },
    @{ Row = 857; A = "Other"; C = "1820"; D = "LAS"; E = "C"; F = 'Turn off wireless microphones by pressing "MUTE" button on mics.'; Height = 30 },
    @{ Row = 858; A = "Other"; C = "1820"; D = "LAS"; E = "C"; F = 'Turn off wireless microphone receivers by pressing "POWER" button and not "SYNC" button. '; Height = 30 },
    @{ Row = 859; A = "Other"; C = "1820"; D = "LAS"; E = "C"; F = "PLEAS BE ON TIME - GO EARLY - GUEST PROF ENDING EARLY TODAY ANYTIME FROM 18:20 - 18:30 pm. "; Height = 30 }
)

$logDate = 42699

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("A$row").Value = $r.A
    $ws.Range("B$row").Value = $logDate
    $ws.Range("C$row").Value = $r.C
    $ws.Range("D$row").Value = $r.D
    $ws.Range("E$row").Value = $r.E
    $ws.Range("F$row").Value = $r.F
    if ($r.Height -gt 0) {
        $ws.Rows.Item($row).RowHeight = $r.Height
    }
}

# Row 859's note is emphasised (bold) in the source workbook, matching the
# other "important reminder" notes elsewhere in column F.
$ws.Range("F859").Font.Bold = $true

[void]$ws.Range("A859").Select()
